$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Login" sheet: drop the three trailing blank rows (6-8) that used to
# pad the sheet and leave the selection sitting on the now-empty area
# right below the data.
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Rows("6:8").Delete() | Out-Null
$loginSheet.Range("A6:A9").EntireRow.Select() | Out-Null

# ---------------------------------------------------------------------
# New "Jobs" sheet, inserted right after "Login", mirroring its
# Run / TCName / Priority layout and formatting.
# ---------------------------------------------------------------------
$jobsSheet = $wb.Worksheets.Add($null, $loginSheet)
$jobsSheet.Name = "Jobs"

$jobsSheet.Cells.Item(1, 1).Value = "Run"
$jobsSheet.Cells.Item(1, 2).Value = "TCName"
$jobsSheet.Cells.Item(1, 3).Value = "Priority"

$jobsSheet.Cells.Item(2, 1).Value = "ON"
$jobsSheet.Cells.Item(2, 2).Value = "SafeWay_Jobs1"
$jobsSheet.Cells.Item(2, 3).Value = 1

$jobsSheet.Cells.Item(5, 1).Value = "ON"
$jobsSheet.Cells.Item(5, 2).Value = "SafeWay_Jobs4"
$jobsSheet.Cells.Item(5, 3).Value = 1

$jobsSheet.Cells.Item(4, 1).Value = "ON"
$jobsSheet.Cells.Item(4, 2).Value = "SafeWay_Jobs3"
$jobsSheet.Cells.Item(4, 3).Value = 1

$jobsSheet.Cells.Item(3, 1).Value = "ON"
$jobsSheet.Cells.Item(3, 2).Value = "SafeWay_Jobs2"
$jobsSheet.Cells.Item(3, 3).Value = 1

# Reuse the same header/body cell styles as "Login".
$loginSheet.Range("A1:C5").Copy()
$jobsSheet.Range("A1:C5").PasteSpecial(-4122)

# Match column B's auto-fit width.
$jobsSheet.Columns.Item(2).ColumnWidth = 51.67

$jobsSheet.Range("B4").Select() | Out-Null
$jobsSheet.Activate() | Out-Null
